# Generate Report for Handback
# Updates the handoff/handback timestamps recorded for the
# "dea62695-3c19-4b4d-8586-a527fcf8cb49.md" file (row 3 in every sheet)
# after a fresh handback Xliff was generated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-01 12:54:26"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-01 12:54:21"
$wsZhCn.Range("K3").Value = "2016-09-01 12:54:38"

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-01 12:54:26"
$wsDeDe.Range("K3").Value = "2016-09-01 12:54:45"
